# Se modifica el campo fecha_pago de las rentas fijas en el informe Nevasa
# Ahora utiliza el siguiente dia habil como fecha de pago (solo fines de semana por ahora).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-6 correspond to the "RENTA FIJA" operations whose fecha_pago (column B)
# must move forward to the next business day: 2024-10-08 (45573) -> 2024-10-09 (45574)
foreach ($row in 3..6) {
    $ws.Cells.Item($row, 2).Value = 45574
}
